# Update countries & provincias Spain
# Refresh of COVID country data: updated totals for several countries and
# updated "last refreshed" timestamp. A handful of adjacent country rows
# also swap places because the refreshed totals changed their sort order
# (the sheet is sorted by column B, "Casos totales", descending).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (row 1) ---
$ws.Range("A1").Value2 = "Datos actualizados a 23 de Junio de 2020 a las 23:41"

# --- Estados Unidos (row 4): updated totals, same rank ---
$ws.Range("B4").Value2 = 2420139
$ws.Range("C4").Value2 = 31986
$ws.Range("D4").Value2 = 1008535
$ws.Range("E4").Value2 = 1288193
$ws.Range("G4").Value2 = 801
$ws.Range("H4").Value2 = 123411

# --- Brasil (row 5): updated totals, same rank ---
$ws.Range("B5").Value2 = 1145906
$ws.Range("C5").Value2 = 34558
$ws.Range("E5").Value2 = 499153
$ws.Range("G5").Value2 = 1242
$ws.Range("H5").Value2 = 52649

# --- Peru (row 10): updated totals, same rank ---
$ws.Range("B10").Value2 = 260810
$ws.Range("C10").Value2 = 3363
$ws.Range("D10").Value2 = 148437
$ws.Range("E10").Value2 = 103969
$ws.Range("G10").Value2 = 181
$ws.Range("H10").Value2 = 8404

# --- Alemania (row 14): updated totals, same rank ---
$ws.Range("B14").Value2 = 192665
$ws.Range("C14").Value2 = 546
$ws.Range("E14").Value2 = 7979
$ws.Range("G14").Value2 = 17
$ws.Range("H14").Value2 = 8986

# --- Barein (row 50): updated totals, same rank ---
$ws.Range("B50").Value2 = 23062
$ws.Range("C50").Value2 = 655
$ws.Range("E50").Value2 = 6133

# --- Uganda / Ruanda (rows 140-141) swap rank: Ruanda's refreshed totals
#     now outrank Uganda, so Ruanda moves up to row 140 (with new data)
#     and Uganda (unchanged data) drops to row 141 ---
$ws.Range("A140").Value2 = "Ruanda"
$ws.Range("B140").Value2 = 798
$ws.Range("C140").Value2 = 11
$ws.Range("D140").Value2 = 371
$ws.Range("E140").Value2 = 425
$ws.Range("H140").Value2 = 2

$ws.Range("A141").Value2 = "Uganda"
$ws.Range("B141").Value2 = 797
$ws.Range("C141").Value2 = 23
$ws.Range("D141").Value2 = 699
$ws.Range("E141").Value2 = 98
$ws.Range("H141").Value2 = 0

# --- Fiyi / Dominica (rows 202-203) swap rank (totals tied, order flips) ---
$ws.Range("A202").Value2 = "Dominica"
$ws.Range("A203").Value2 = "Fiyi"

# --- Islas Malvinas / Groenlandia (rows 208-209) swap rank (totals tied) ---
$ws.Range("A208").Value2 = "Groenlandia"
$ws.Range("A209").Value2 = "Islas Malvinas"

# --- Montserrat / Seychelles (rows 211-212) swap rank, data trades places ---
$ws.Range("A211").Value2 = "Seychelles"
$ws.Range("D211").Value2 = 11
$ws.Range("H211").Value2 = 0

$ws.Range("A212").Value2 = "Montserrat"
$ws.Range("D212").Value2 = 10
$ws.Range("H212").Value2 = 1
